$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $d, $e) {
    $cellD = $ws.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

# Row 2: Bitcoin
Set-Row 2 "22.108.29" "  -1.57%  "
# Row 3: Ethereum
Set-Row 3 "1.554.34" "  -1.12%  "
# Row 4: TetherUSD
Set-Row 4 "0.9995" "  -0.19%  "
# Row 5: USDC
Set-Row 5 "1.0000" "  -0.14%  "
# Row 6: BNB
Set-Row 6 "287.83" "  -0.19%  "
# Row 7: XRP
Set-Row 7 "0.3806" "  +2.61%  "
# Row 8: Cardano
Set-Row 8 "0.3289" "  -0.69%  "
# Row 9: OKB
Set-Row 9 "43.12" "  -10.80%  "
# Row 10: Polygon
Set-Row 10 "1.138" "  +0.37%  "
# Row 11: Dogecoin
Set-Row 11 "0.07342" "  -2.01%  "
# Row 12: BinanceUSD
Set-Row 12 "0.9996" "  -0.21%  "
# Row 13: Solana
Set-Row 13 "20.16" "  -2.75%  "
# Row 14: Polkadot
Set-Row 14 "5.811" "  -1.98%  "
# Row 15: Chainlink
Set-Row 15 "6.800" "  -1.19%  "
# Row 16: WrappedEther
Set-Row 16 "1.554.76" "  -1.06%  "
# Row 17: ShibaInu
Set-Row 17 "0.00001094" "  -2.22%  "
# Row 18: TRON
Set-Row 18 "0.06621" "  -1.82%  "
# Row 19: Litecoin
Set-Row 19 "85.79" "  -1.99%  "

# Row 20 and 21 swap (Dai <-> Uniswap)
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-Row 20 "6.397" "  +0.75%  "

$ws.Cells.Item(21, 2).Value = "Dai"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-Row 21 "0.9992" "  -0.19%  "

# Row 22: Avalanche
Set-Row 22 "16.09" "  -2.63%  "
# Row 23: Cosmos
Set-Row 23 "11.69" "  -2.99%  "
# Row 24: WrappedBTC
Set-Row 24 "22.093.55" "  -1.59%  "
# Row 25: Toncoin
Set-Row 25 "2.318" "  -2.68%  "
# Row 26: LidoDAOToken
Set-Row 26 "2.521" "  -2.19%  "
# Row 27: Monero
Set-Row 27 "150.39" "  -2.10%  "
# Row 28: EthereumClassic
Set-Row 28 "19.12" "  -2.93%  "
# Row 29: HuobiToken
Set-Row 29 "4.912" "  -2.25%  "
# Row 30: BitcoinCash
Set-Row 30 "121.43" "  -2.50%  "
# Row 31: WrappedliquidstakedEther2.0
Set-Row 31 "1.724.63" "  -1.28%  "
# Row 32: ImmutableX
Set-Row 32 "1.072" "  +1.07%  "
# Row 33: Filecoin
Set-Row 33 "5.925" "  -3.17%  "
# Row 34: WEMIXTOKEN
Set-Row 34 "1.866" "  -7.34%  "
# Row 35: Stellar
Set-Row 35 "0.08233" "  -1.42%  "
# Row 36: FraxShare
Set-Row 36 "9.300" "  -5.00%  "
# Row 37: VeChain
Set-Row 37 "0.02326" "  -5.55%  "
# Row 38: InternetComputer(DFINITY)
Set-Row 38 "5.285" "  -1.01%  "
# Row 39: Hedera
Set-Row 39 "0.06233" "  -2.78%  "
# Row 40: Algorand
Set-Row 40 "0.2162" "  -4.67%  "
# Row 41: TrustWalletToken
Set-Row 41 "1.249" "  -3.17%  "
# Row 42: Aptos
Set-Row 42 "11.04" "  -2.28%  "

# Row 43 and 44 swap (Frax <-> TheSandbox)
$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-Row 43 "0.6035" "  -4.65%  "

$ws.Cells.Item(44, 2).Value = "Frax"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-Row 44 "0.9997" "  -0.14%  "

# Row 45: EnergySwap
Set-Row 45 "13.64" "  -1.99%  "
# Row 46: PancakeSwap (only E changes)
$ws.Cells.Item(46, 5).Value = "  -1.10%  "
# Row 47: Decentraland
Set-Row 47 "0.5842" "  -5.32%  "
# Row 48: NEARProtocol
Set-Row 48 "1.986" "  -3.53%  "
# Row 49: Quant
Set-Row 49 "122.07" "  -2.98%  "
# Row 50: EOS
Set-Row 50 "1.174" "  -3.07%  "
# Row 51: Cronos (only D changes)
$cellD51 = $ws.Cells.Item(51, 4)
$cellD51.NumberFormat = "@"
$cellD51.Value = "0.07008"
